# Weekly data update: insert two new rows of "Pera" price data for
# Feria Lagunitas de Puerto Montt (Packham's Triumph, Primera & Segunda)
# right above the existing row 157, shifting the rest of the table down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 157 (existing rows 157:165 shift down to 159:167).
$ws.Range("A157:A158").EntireRow.Insert()

# New row 157: Packham's Triumph, Primera
$ws.Range("A157").Value = 4
$ws.Range("B157").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C157").Value = "Los Lagos"
$ws.Range("D157").Value = 44509
$ws.Range("E157").Value = 10
$ws.Range("F157").Value = "Fruta"
$ws.Range("G157").Value = 100104
$ws.Range("H157").Value = "Frutos de pepita"
$ws.Range("I157").Value = 100104005
$ws.Range("J157").Value = "Pera"
$ws.Range("K157").Value = "Packham's Triumph"
$ws.Range("L157").Value = "Primera"
$ws.Range("M157").Value = 500
$ws.Range("N157").Value = 15000
$ws.Range("O157").Value = 16000
$ws.Range("P157").Value = 15500
$ws.Range("Q157").Value = '$/caja 15 kilos empedrada'
$ws.Range("R157").Value = "Región de O'Higgins"
$ws.Range("S157").Value = 1033
$ws.Range("T157").Value = 15

# New row 158: Packham's Triumph, Segunda
$ws.Range("A158").Value = 4
$ws.Range("B158").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C158").Value = "Los Lagos"
$ws.Range("D158").Value = 44509
$ws.Range("E158").Value = 10
$ws.Range("F158").Value = "Fruta"
$ws.Range("G158").Value = 100104
$ws.Range("H158").Value = "Frutos de pepita"
$ws.Range("I158").Value = 100104005
$ws.Range("J158").Value = "Pera"
$ws.Range("K158").Value = "Packham's Triumph"
$ws.Range("L158").Value = "Segunda"
$ws.Range("M158").Value = 300
$ws.Range("N158").Value = 13000
$ws.Range("O158").Value = 13000
$ws.Range("P158").Value = 13000
$ws.Range("Q158").Value = '$/caja 15 kilos empedrada'
$ws.Range("R158").Value = "Región de O'Higgins"
$ws.Range("S158").Value = 867
$ws.Range("T158").Value = 15
